$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "SCREE" + "N" + "ING REPORT" (3 runs, with a _GoBack bookmark sitting
#    between run 1 and run 3) collapse into a single run reading
#    "SCREEING REPORT" - this also removes the _GoBack bookmark from here.
#    Use InsertXML on the exact text range so the surviving run keeps the
#    first run's rsidRPr attribute instead of getting a "fresh" run.
# ---------------------------------------------------------------------------
$titleRange = $d.Content
$titleFound = $titleRange.Find.Execute("SCREENING REPORT")

$titleXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body><w:p>' +
  '<w:r w:rsidRPr="005219F9"><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>SCREEING REPORT</w:t></w:r>' +
  '</w:p></w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

if ($titleFound) {
    $titleRange.InsertXML($titleXml) | Out-Null
}

# Belt-and-braces: the InsertXML above spans the old bookmark's collapsed
# position, so _GoBack is already gone from the title - but make sure, in
# case of a quirky re-seat, it isn't left stranded in the title paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $staleBm = $d.Bookmarks.Item("_GoBack")
    if ($staleBm.Start -lt $titleRange.End) {
        $staleBm.Delete()
    }
}

# ---------------------------------------------------------------------------
# 2. Date text update. Use InsertXML (instead of Find/Replace text) so the
#    surviving run keeps its original w:rsidR="00D75F60" and full <w:rPr>
#    instead of being replaced by a "fresh" unattributed run.
# ---------------------------------------------------------------------------
$dateRange = $d.Content
$dateFound = $dateRange.Find.Execute("2021-03-04")

$dateXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body><w:p>' +
  '<w:r w:rsidR="00D75F60"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>2021-02-26</w:t></w:r>' +
  '</w:p></w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

if ($dateFound) {
    $dateRange2 = $d.Range($dateRange.Start, $dateRange.End)
    $dateRange2.InsertXML($dateXml) | Out-Null
}

# ---------------------------------------------------------------------------
# 3. Re-home the _GoBack bookmark: it now needs to sit right after the
#    trailing page-break run, inside the very last paragraph of the body
#    (immediately before that paragraph's end mark).
#    Bookmarks.Add() on an already-collapsed Range mis-fires in this host,
#    so insert a throw-away character, bookmark the 1-char span, then
#    delete the character again - the bookmark collapses cleanly in place.
# ---------------------------------------------------------------------------
$endRange = $d.Range($d.Content.End, $d.Content.End)
$endRange.MoveEnd(1, -1) | Out-Null            # step back over the final paragraph mark
$dropPoint = $endRange.End                      # right after the page break run

$marker = $d.Range($dropPoint, $dropPoint)
$marker.InsertAfter("X")
$markerSpan = $d.Range($dropPoint, $dropPoint + 1)
$d.Bookmarks.Add("_GoBack", $markerSpan)
$markerSpan2 = $d.Range($dropPoint, $dropPoint + 1)
$markerSpan2.Delete()
